$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '22.350.10'
$ws.Range('E2').Value = '  -4.84%  '

$ws.Range('D3').Value = '1.570.13'
$ws.Range('E3').Value = '  -4.84%  '

$ws.Range('E4').Value = '  -0.05%  '

$ws.Range('E5').Value = '  +0.00%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range('D6').Value = '289.86'
$ws.Range('E6').Value = '  -3.44%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range('D7').Value = '0.3756'
$ws.Range('E7').Value = '  -0.73%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range('D8').Value = '49.57'
$ws.Range('E8').Value = '  -2.65%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range('D9').Value = '0.3407'
$ws.Range('E9').Value = '  -2.64%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range('D10').Value = '1.168'
$ws.Range('E10').Value = '  -4.89%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range('D11').Value = '0.07677'
$ws.Range('E11').Value = '  -4.75%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range('D12').Value = '1.001'
$ws.Range('E12').Value = '  +0.01%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range('D13').Value = '21.39'
$ws.Range('E13').Value = '  -2.96%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range('D14').Value = '6.031'
$ws.Range('E14').Value = '  -4.60%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range('D15').Value = '6.950'
$ws.Range('E15').Value = '  -4.47%  '

$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '1.577.78'
$ws.Range('E16').Value = '  -4.12%  '

$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").NumberFormat = "@"
$ws.Range('D17').Value = '0.00001138'
$ws.Range('E17').Value = '  -6.01%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range('D18').Value = '90.32'
$ws.Range('E18').Value = '  -5.19%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range('D19').Value = '0.06725'
$ws.Range('E19').Value = '  -3.64%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range('D20').Value = '1.001'
$ws.Range('E20').Value = '  +0.02%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range('D21').Value = '6.254'
$ws.Range('E21').Value = '  -5.89%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range('D22').Value = '16.65'
$ws.Range('E22').Value = '  -4.78%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range('D23').Value = '0.5300'
$ws.Range('E23').Value = '  -8.10%  '

$ws.Range('E24').Value = '  -4.13%  '

$ws.Range('D25').Value = '22.336.39'
$ws.Range('E25').Value = '  -4.89%  '

$ws.Range('E26').Value = '  -1.25%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range('D27').Value = '2.793'
$ws.Range('E27').Value = '  -7.55%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range('D28').Value = '20.17'
$ws.Range('E28').Value = '  -4.54%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range('D29').Value = '145.29'
$ws.Range('E29').Value = '  -3.94%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range('D30').Value = '4.993'
$ws.Range('E30').Value = '  -3.46%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range('D31').Value = '125.80'
$ws.Range('E31').Value = '  -4.42%  '

$ws.Range('D32').Value = '1.743.58'
$ws.Range('E32').Value = '  -4.98%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range('D33').Value = '1.023'
$ws.Range('E33').Value = '  +3.28%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range('D34').Value = '6.232'
$ws.Range('E34').Value = '  -9.58%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range('D35').Value = '2.014'
$ws.Range('E35').Value = '  -6.30%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range('D36').Value = '10.14'
$ws.Range('E36').Value = '  -8.91%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range('D37').Value = '0.08524'
$ws.Range('E37').Value = '  -3.01%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range('D38').Value = '0.02541'
$ws.Range('E38').Value = '  -7.02%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range('D39').Value = '0.2328'
$ws.Range('E39').Value = '  -4.01%  '

$ws.Range('E40').Value = '  -6.50%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range('D41').Value = '1.328'
$ws.Range('E41').Value = '  +2.44%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range('D42').Value = '0.06411'
$ws.Range('E42').Value = '  -7.00%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range('D43').Value = '11.79'
$ws.Range('E43').Value = '  -8.36%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range('D44').Value = '0.6423'

$ws.Range("D45").NumberFormat = "@"
$ws.Range('D45').Value = '14.23'
$ws.Range('E45').Value = '  -8.63%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range('D46').Value = '0.9998'
$ws.Range('E46').Value = '  +0.02%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range('D47').Value = '0.5998'
$ws.Range('E47').Value = '  -6.15%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range('D48').Value = '3.759'
$ws.Range('E48').Value = '  -4.05%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range('D49').Value = '2.100'
$ws.Range('E49').Value = '  -6.93%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range('D50').Value = '1.278'
$ws.Range('E50').Value = '  +3.23%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range('D51').Value = '124.79'
$ws.Range('E51').Value = '  -1.72%  '
